# Delete the rows for years 2004-2009 (rows 2 through 7). Excel's
# Range.Delete shifts the remaining rows (2010-2019, previously rows 8-17)
# up so they become rows 2-11, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:I7").EntireRow.Delete()
